$d = $word.ActiveDocument

# Update the date/weekday heading.
$d.Content.Find.Execute("2023-10-20 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-10-21 Saturday", 2)

# Update the division problems in the first (and only) table.
# Cells are addressed directly by (row, column) so that overlapping
# old/new values (e.g. "19÷3=" and "67÷7=") never get double-replaced.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "80÷3="
$t.Cell(1, 2).Range.Text = "76÷2="
$t.Cell(1, 3).Range.Text = "57÷4="
$t.Cell(1, 4).Range.Text = "77÷6="
$t.Cell(1, 5).Range.Text = "67÷7="

$t.Cell(5, 1).Range.Text = "67÷7="
$t.Cell(5, 2).Range.Text = "59÷9="
$t.Cell(5, 3).Range.Text = "82÷8="
$t.Cell(5, 4).Range.Text = "31÷5="
$t.Cell(5, 5).Range.Text = "19÷3="

$t.Cell(9, 1).Range.Text = "64÷9="
$t.Cell(9, 2).Range.Text = "92÷4="
$t.Cell(9, 3).Range.Text = "32÷4="
$t.Cell(9, 4).Range.Text = "96÷3="
$t.Cell(9, 5).Range.Text = "81÷7="

$t.Cell(13, 1).Range.Text = "78÷6="
$t.Cell(13, 2).Range.Text = "68÷6="
$t.Cell(13, 3).Range.Text = "37÷2="
$t.Cell(13, 4).Range.Text = "56÷4="
$t.Cell(13, 5).Range.Text = "94÷4="

$t.Cell(17, 1).Range.Text = "15÷7="
$t.Cell(17, 2).Range.Text = "81÷7="
$t.Cell(17, 3).Range.Text = "77÷3="
$t.Cell(17, 4).Range.Text = "78÷7="
$t.Cell(17, 5).Range.Text = "16÷8="
